$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds an Excel date-serial value for every data row
# (rows 2-433). The daily automated update bumps that date by one day
# (serial 46074 -> 46075) for every row.
$lastRow = 433
$ws.Range("C2:C$lastRow").Value = 46075
